$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "CÔ DIỄM" (1st tab) — insert a new dated row at row 25, and
# append two new journal rows (Duy lends cô Diễm 5tr then the standard
# -10000 "vay 10tr" line) at the bottom of the ledger.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Insert a new row above the current row 25; Excel shifts every
# formula/reference below it down by one automatically (shared
# formulas, the J3 SUM range, etc.)
$ws1.Rows("25:25").Insert()

# The freshly inserted row inherits a slightly different auto style;
# re-copy the real formatting from the row right below (which is the
# old row 25, now shifted to row 26) so H25/I25 match the rest of the
# column.
$ws1.Cells.Item(26, 8).Copy($ws1.Cells.Item(25, 8))
$ws1.Cells.Item(26, 9).Copy($ws1.Cells.Item(25, 9))

$ws1.Cells.Item(25, 8).Value = "31/05/2023"
$ws1.Cells.Item(25, 9).Value = 10

# Append the two new ledger entries after the last existing row (old
# row 112, now row 113).
$ws1.Cells.Item(114, 2).Value = "Duy chyển cho cô Diễm 5tr"
$ws1.Cells.Item(114, 3).Value = 5000
$ws1.Cells.Item(114, 4).Formula = "=D113+C114"

$ws1.Cells.Item(115, 2).Value = "Duy cho cô Diễm vay 10tr"
$ws1.Cells.Item(115, 3).Value = -10000
$ws1.Cells.Item(115, 4).Formula = "=D114+C115"

# Restore the matching cell formatting (border/number format/font) on
# the new rows by copying from the row above.
$ws1.Cells.Item(113, 2).Copy($ws1.Cells.Item(114, 2))
$ws1.Cells.Item(113, 3).Copy($ws1.Cells.Item(114, 3))
$ws1.Cells.Item(113, 4).Copy($ws1.Cells.Item(114, 4))
$ws1.Cells.Item(114, 2).Copy($ws1.Cells.Item(115, 2))
$ws1.Cells.Item(114, 3).Copy($ws1.Cells.Item(115, 3))
$ws1.Cells.Item(114, 4).Copy($ws1.Cells.Item(115, 4))

$ws1.Cells.Item(114, 2).Value = "Duy chyển cho cô Diễm 5tr"
$ws1.Cells.Item(114, 3).Value = 5000
$ws1.Cells.Item(114, 4).Formula = "=D113+C114"

$ws1.Cells.Item(115, 2).Value = "Duy cho cô Diễm vay 10tr"
$ws1.Cells.Item(115, 3).Value = -10000
$ws1.Cells.Item(115, 4).Formula = "=D114+C115"

# Scroll position / selection left behind by the author after the edit.
$ws1.Application.Goto($ws1.Range("I26"))
$ws1.Range("I26").Select()

# ---------------------------------------------------------------------
# Sheet "Sheet1" (6th tab) — new debt breakdown table in A1:B4 feeding
# the existing totals in row 15-17.
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

$ws6.Cells.Item(1, 1).Value = "tran ngoc quang"
$ws6.Cells.Item(1, 2).Value = 400
$ws6.Cells.Item(2, 1).Value = "nguyen khanh duy"
$ws6.Cells.Item(2, 2).Value = 1000
$ws6.Cells.Item(3, 1).Value = "co diem no"
$ws6.Cells.Item(3, 2).Value = 700
$ws6.Cells.Item(4, 1).Value = "vo thi dang nga"
$ws6.Cells.Item(4, 2).Value = 400

$ws6.Cells.Item(15, 2).Formula = "=SUM(B1:B13)"
$ws6.Cells.Item(16, 2).Value = 5000

$ws6.Range("B5").Select()
